$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 6 (entire row) - data shrinks from 6 rows to 5 rows
$ws.Rows.Item(6).Delete()

# Update column widths (character-width units; Excel COM ColumnWidth
# differs from the stored <col width> by a fixed 0.8333... padding,
# so subtract that offset to land on the exact target width.)
$ws.Columns.Item(12).ColumnWidth = 6.166666666666666  # L: 8 -> 7
$ws.Columns.Item(20).ColumnWidth = 8.166666666666666  # T: 8 -> 9
$ws.Columns.Item(30).ColumnWidth = 6.166666666666666  # AD: 8 -> 7
$ws.Columns.Item(34).ColumnWidth = 6.166666666666666  # AH: 8 -> 7

# Replace data rows 2-5 with the new dataset

# Row 2
$ws.Range("A2").Value = 45058.50694444445
$ws.Range("B2").Value = 10.726
$ws.Range("C2").Value = 7.333
$ws.Range("D2").Value = 3.404
$ws.Range("E2").Value = 23.56
$ws.Range("F2").Value = 17.15
$ws.Range("G2").Value = 8.176
$ws.Range("H2").Value = 24.228
$ws.Range("I2").Value = 13.347
$ws.Range("J2").Value = 5.245
$ws.Range("K2").Value = 7.323
$ws.Range("L2").Value = 9.308999999999999
$ws.Range("M2").Value = 10.191
$ws.Range("N2").Value = 2.44
$ws.Range("O2").Value = 8.647
$ws.Range("P2").Value = 11.655
$ws.Range("Q2").Value = 7.955
$ws.Range("R2").Value = 2.648
$ws.Range("S2").Value = 1.093
$ws.Range("T2").Value = 124.223
$ws.Range("U2").Value = 23.834
$ws.Range("V2").Value = 7.982
$ws.Range("W2").Value = 14.964
$ws.Range("X2").Value = 8.048999999999999
$ws.Range("Y2").Value = 2.19
$ws.Range("Z2").Value = 13.597
$ws.Range("AA2").Value = 7.05
$ws.Range("AB2").Value = 6.629
$ws.Range("AC2").Value = 7.562
$ws.Range("AD2").Value = 9.942
$ws.Range("AE2").Value = 2.682
$ws.Range("AF2").Value = 21.628
$ws.Range("AG2").Value = 4.075
$ws.Range("AH2").Value = 9.978

# Row 3
$ws.Range("A3").Value = 45058.51388888889
$ws.Range("B3").Value = 1.224
$ws.Range("C3").Value = 0.649
$ws.Range("D3").Value = 1.25
$ws.Range("E3").Value = 3.047
$ws.Range("F3").Value = 1.289
$ws.Range("G3").Value = 0.845
$ws.Range("H3").Value = 8.363
$ws.Range("I3").Value = 1.714
$ws.Range("J3").Value = 0.641
$ws.Range("K3").Value = 0.245
$ws.Range("L3").Value = 1.184
$ws.Range("M3").Value = 1.473
$ws.Range("N3").Value = 0.17
$ws.Range("O3").Value = 1.128
$ws.Range("P3").Value = 1.444
$ws.Range("Q3").Value = 1.472
$ws.Range("R3").Value = 1.21
$ws.Range("S3").Value = 0.326
$ws.Range("T3").Value = 9.956
$ws.Range("U3").Value = 3.578
$ws.Range("V3").Value = 1.041
$ws.Range("W3").Value = 1.977
$ws.Range("X3").Value = 0.955
$ws.Range("Y3").Value = 0.5570000000000001
$ws.Range("Z3").Value = 4.471
$ws.Range("AA3").Value = 0.92
$ws.Range("AB3").Value = 1.094
$ws.Range("AC3").Value = 1.203
$ws.Range("AD3").Value = 1.264
$ws.Range("AE3").Value = 1.136
$ws.Range("AF3").Value = 8.289
$ws.Range("AG3").Value = 0.316
$ws.Range("AH3").Value = 1.312

# Row 4
$ws.Range("A4").Value = 45058.52083333334
$ws.Range("B4").Value = 8.964
$ws.Range("C4").Value = 6.585
$ws.Range("D4").Value = 1.055
$ws.Range("E4").Value = 19.821
$ws.Range("F4").Value = 15.551
$ws.Range("G4").Value = 6.978
$ws.Range("H4").Value = 24.499
$ws.Range("I4").Value = 11.023
$ws.Range("J4").Value = 4.788
$ws.Range("K4").Value = 6.769
$ws.Range("L4").Value = 7.937
$ws.Range("M4").Value = 8.561999999999999
$ws.Range("N4").Value = 2.124
$ws.Range("O4").Value = 7.143
$ws.Range("P4").Value = 9.929
$ws.Range("Q4").Value = 6.304
$ws.Range("R4").Value = 0.874
$ws.Range("S4").Value = 0.431
$ws.Range("T4").Value = 101.353
$ws.Range("U4").Value = 19.79
$ws.Range("V4").Value = 6.594
$ws.Range("W4").Value = 12.96
$ws.Range("X4").Value = 6.875
$ws.Range("Y4").Value = 1.19
$ws.Range("Z4").Value = 12.411
$ws.Range("AA4").Value = 5.824
$ws.Range("AB4").Value = 5.303
$ws.Range("AC4").Value = 6.201
$ws.Range("AD4").Value = 8.351000000000001
$ws.Range("AE4").Value = 0.722
$ws.Range("AF4").Value = 21.971
$ws.Range("AG4").Value = 3.565
$ws.Range("AH4").Value = 8.242000000000001

# Row 5
$ws.Range("A5").Value = 45058.52777777778
$ws.Range("B5").Value = 14.28
$ws.Range("C5").Value = 10.63
$ws.Range("D5").Value = 1.03
$ws.Range("E5").Value = 31.33
$ws.Range("F5").Value = 25.27
$ws.Range("G5").Value = 11.18
$ws.Range("H5").Value = 41.98
$ws.Range("I5").Value = 17.42
$ws.Range("J5").Value = 7.75
$ws.Range("K5").Value = 11.21
$ws.Range("L5").Value = 12.56
$ws.Range("M5").Value = 13.43
$ws.Range("N5").Value = 3.49
$ws.Range("O5").Value = 11.28
$ws.Range("P5").Value = 15.94
$ws.Range("Q5").Value = 9.640000000000001
$ws.Range("R5").Value = 0.72
$ws.Range("S5").Value = 0.53
$ws.Range("T5").Value = 164.34
$ws.Range("U5").Value = 31.49
$ws.Range("V5").Value = 10.41
$ws.Range("W5").Value = 21.02
$ws.Range("X5").Value = 11.06
$ws.Range("Y5").Value = 1.66
$ws.Range("Z5").Value = 20.97
$ws.Range("AA5").Value = 9.199999999999999
$ws.Range("AB5").Value = 8.220000000000001
$ws.Range("AC5").Value = 9.640000000000001
$ws.Range("AD5").Value = 13.21
$ws.Range("AE5").Value = 0.53
$ws.Range("AF5").Value = 38.11
$ws.Range("AG5").Value = 5.79
$ws.Range("AH5").Value = 13.02
